$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A2").Value = "Última actualización: 30/12/2025 15:58:16"
$ws1.Range("A3").Value = "Total filas: 361"

$ws1.Cells.Item(340, 2).Value = "15:58:05"
$ws1.Cells.Item(340, 3).Value = "16:00"
$ws1.Cells.Item(340, 4).Value = "10_OLMOS"
$ws1.Cells.Item(340, 5).Value = 2
$ws1.Cells.Item(340, 6).Value = "LP1912"
$ws1.Cells.Item(340, 7).Value = "30/12/2025"

$ws1.Cells.Item(341, 2).Value = "15:58:05"
$ws1.Cells.Item(341, 3).Value = "16:03"
$ws1.Cells.Item(341, 4).Value = "16_SANTA ANA"
$ws1.Cells.Item(341, 5).Value = 5
$ws1.Cells.Item(341, 6).Value = "LP1912"
$ws1.Cells.Item(341, 7).Value = "30/12/2025"

$ws1.Cells.Item(342, 2).Value = "15:58:05"
$ws1.Cells.Item(342, 3).Value = "16:05"
$ws1.Cells.Item(342, 4).Value = "23_HERNANDEZ"
$ws1.Cells.Item(342, 5).Value = 7
$ws1.Cells.Item(342, 6).Value = "LP1912"
$ws1.Cells.Item(342, 7).Value = "30/12/2025"

$ws1.Cells.Item(343, 2).Value = "15:58:05"
$ws1.Cells.Item(343, 3).Value = "16:11"
$ws1.Cells.Item(343, 4).Value = "16_SANTA ANA"
$ws1.Cells.Item(343, 5).Value = 13
$ws1.Cells.Item(343, 6).Value = "LP1912"
$ws1.Cells.Item(343, 7).Value = "30/12/2025"

$ws1.Cells.Item(344, 2).Value = "15:58:05"
$ws1.Cells.Item(344, 3).Value = "16:20"
$ws1.Cells.Item(344, 4).Value = "215C_EL PATO"
$ws1.Cells.Item(344, 5).Value = 22
$ws1.Cells.Item(344, 6).Value = "LP1912"
$ws1.Cells.Item(344, 7).Value = "30/12/2025"

$ws1.Cells.Item(345, 2).Value = "15:58:05"
$ws1.Cells.Item(345, 3).Value = "16:21"
$ws1.Cells.Item(345, 4).Value = "26_HERNANDEZ"
$ws1.Cells.Item(345, 5).Value = 23
$ws1.Cells.Item(345, 6).Value = "LP1912"
$ws1.Cells.Item(345, 7).Value = "30/12/2025"

$ws1.Cells.Item(346, 2).Value = "15:58:05"
$ws1.Cells.Item(346, 3).Value = "16:27"
$ws1.Cells.Item(346, 4).Value = "16_SANTA ANA"
$ws1.Cells.Item(346, 5).Value = 29
$ws1.Cells.Item(346, 6).Value = "LP1912"
$ws1.Cells.Item(346, 7).Value = "30/12/2025"

$ws1.Cells.Item(347, 2).Value = "15:58:05"
$ws1.Cells.Item(347, 3).Value = "16:29"
$ws1.Cells.Item(347, 4).Value = "10_OLMOS"
$ws1.Cells.Item(347, 5).Value = 31
$ws1.Cells.Item(347, 6).Value = "LP1912"
$ws1.Cells.Item(347, 7).Value = "30/12/2025"

$ws1.Cells.Item(348, 2).Value = "15:58:05"
$ws1.Cells.Item(348, 3).Value = "16:35"
$ws1.Cells.Item(348, 4).Value = "23_HERNANDEZ"
$ws1.Cells.Item(348, 5).Value = 37
$ws1.Cells.Item(348, 6).Value = "LP1912"
$ws1.Cells.Item(348, 7).Value = "30/12/2025"

$ws1.Cells.Item(349, 2).Value = "15:58:05"
$ws1.Cells.Item(349, 3).Value = "16:37"
$ws1.Cells.Item(349, 4).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(349, 5).Value = 39
$ws1.Cells.Item(349, 6).Value = "LP1912"
$ws1.Cells.Item(349, 7).Value = "30/12/2025"

$ws1.Cells.Item(350, 2).Value = "15:58:05"
$ws1.Cells.Item(350, 3).Value = "16:43"
$ws1.Cells.Item(350, 4).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(350, 5).Value = 45
$ws1.Cells.Item(350, 6).Value = "LP1912"
$ws1.Cells.Item(350, 7).Value = "30/12/2025"

$ws1.Cells.Item(351, 2).Value = "15:58:05"
$ws1.Cells.Item(351, 3).Value = "16:48"
$ws1.Cells.Item(351, 4).Value = "15_ABASTO"
$ws1.Cells.Item(351, 5).Value = 50
$ws1.Cells.Item(351, 6).Value = "LP1912"
$ws1.Cells.Item(351, 7).Value = "30/12/2025"

$ws1.Cells.Item(352, 2).Value = "15:58:05"
$ws1.Cells.Item(352, 3).Value = "16:51"
$ws1.Cells.Item(352, 4).Value = "14_ABASTO"
$ws1.Cells.Item(352, 5).Value = 53
$ws1.Cells.Item(352, 6).Value = "LP1912"
$ws1.Cells.Item(352, 7).Value = "30/12/2025"

$ws1.Cells.Item(353, 2).Value = "15:58:05"
$ws1.Cells.Item(353, 3).Value = "16:56"
$ws1.Cells.Item(353, 4).Value = "17_179 Y 38"
$ws1.Cells.Item(353, 5).Value = 58
$ws1.Cells.Item(353, 6).Value = "LP1912"
$ws1.Cells.Item(353, 7).Value = "30/12/2025"

$ws1.Cells.Item(354, 2).Value = "15:58:05"
$ws1.Cells.Item(354, 3).Value = "16:57"
$ws1.Cells.Item(354, 4).Value = "10_OLMOS"
$ws1.Cells.Item(354, 5).Value = 59
$ws1.Cells.Item(354, 6).Value = "LP1912"
$ws1.Cells.Item(354, 7).Value = "30/12/2025"

$ws1.Cells.Item(355, 2).Value = "15:58:05"
$ws1.Cells.Item(355, 3).Value = "17:05"
$ws1.Cells.Item(355, 4).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(355, 5).Value = 67
$ws1.Cells.Item(355, 6).Value = "LP1912"
$ws1.Cells.Item(355, 7).Value = "30/12/2025"

$ws1.Cells.Item(356, 2).Value = "15:58:05"
$ws1.Cells.Item(356, 3).Value = "17:05"
$ws1.Cells.Item(356, 4).Value = "215A_EL PATO"
$ws1.Cells.Item(356, 5).Value = 67
$ws1.Cells.Item(356, 6).Value = "LP1912"
$ws1.Cells.Item(356, 7).Value = "30/12/2025"

$ws1.Cells.Item(357, 2).Value = "15:58:05"
$ws1.Cells.Item(357, 3).Value = "17:05"
$ws1.Cells.Item(357, 4).Value = "23_HERNANDEZ"
$ws1.Cells.Item(357, 5).Value = 67
$ws1.Cells.Item(357, 6).Value = "LP1912"
$ws1.Cells.Item(357, 7).Value = "30/12/2025"

$ws1.Cells.Item(358, 2).Value = "15:58:05"
$ws1.Cells.Item(358, 3).Value = "17:21"
$ws1.Cells.Item(358, 4).Value = "26_HERNANDEZ"
$ws1.Cells.Item(358, 5).Value = 83
$ws1.Cells.Item(358, 6).Value = "LP1912"
$ws1.Cells.Item(358, 7).Value = "30/12/2025"

$ws1.Cells.Item(359, 2).Value = "15:58:05"
$ws1.Cells.Item(359, 3).Value = "17:24"
$ws1.Cells.Item(359, 4).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(359, 5).Value = 86
$ws1.Cells.Item(359, 6).Value = "LP1912"
$ws1.Cells.Item(359, 7).Value = "30/12/2025"

$ws1.Cells.Item(360, 2).Value = "15:58:05"
$ws1.Cells.Item(360, 3).Value = "17:29"
$ws1.Cells.Item(360, 4).Value = "14_ABASTO"
$ws1.Cells.Item(360, 5).Value = 91
$ws1.Cells.Item(360, 6).Value = "LP1912"
$ws1.Cells.Item(360, 7).Value = "30/12/2025"

$ws1.Cells.Item(361, 2).Value = "15:58:05"
$ws1.Cells.Item(361, 3).Value = "17:31"
$ws1.Cells.Item(361, 4).Value = "15_ABASTO"
$ws1.Cells.Item(361, 5).Value = 93
$ws1.Cells.Item(361, 6).Value = "LP1912"
$ws1.Cells.Item(361, 7).Value = "30/12/2025"

$ws1.Cells.Item(362, 2).Value = "15:58:05"
$ws1.Cells.Item(362, 3).Value = "17:35"
$ws1.Cells.Item(362, 4).Value = "27_EL RETIRO"
$ws1.Cells.Item(362, 5).Value = 97
$ws1.Cells.Item(362, 6).Value = "LP1912"
$ws1.Cells.Item(362, 7).Value = "30/12/2025"

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2").Value = "Última actualización: 30/12/2025 15:58:16"
$ws2.Range("A3").Value = "Total filas: 28"

$ws2.Cells.Item(28, 2).Value = "30/12/2025"
$ws2.Cells.Item(28, 3).Value = "15:58:05"
$ws2.Cells.Item(28, 4).Value = "16:20"
$ws2.Cells.Item(28, 5).Value = "215C_EL PATO"
$ws2.Cells.Item(28, 6).Value = 22
$ws2.Cells.Item(28, 7).Value = "LP1912"

$ws2.Cells.Item(29, 2).Value = "30/12/2025"
$ws2.Cells.Item(29, 3).Value = "15:58:05"
$ws2.Cells.Item(29, 4).Value = "17:05"
$ws2.Cells.Item(29, 5).Value = "215A_EL PATO"
$ws2.Cells.Item(29, 6).Value = 67
$ws2.Cells.Item(29, 7).Value = "LP1912"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A2").Value = "Última actualización: 30/12/2025 15:58:16"
$ws3.Range("A3").Value = "Total filas: 50"

$ws3.Cells.Item(49, 2).Value = "30/12/2025"
$ws3.Cells.Item(49, 3).Value = "15:58:11"
$ws3.Cells.Item(49, 4).Value = "16:14"
$ws3.Cells.Item(49, 5).Value = "215C_LA PLATA"
$ws3.Cells.Item(49, 6).Value = 16
$ws3.Cells.Item(49, 7).Value = "L6203"

$ws3.Cells.Item(50, 2).Value = "30/12/2025"
$ws3.Cells.Item(50, 3).Value = "15:58:16"
$ws3.Cells.Item(50, 4).Value = "16:53"
$ws3.Cells.Item(50, 5).Value = "215B_LP-P MOR-40 Y 115"
$ws3.Cells.Item(50, 6).Value = 55
$ws3.Cells.Item(50, 7).Value = "L6173"

$ws3.Cells.Item(51, 2).Value = "30/12/2025"
$ws3.Cells.Item(51, 3).Value = "15:58:16"
$ws3.Cells.Item(51, 4).Value = "17:21"
$ws3.Cells.Item(51, 5).Value = "215A_LA PLATA"
$ws3.Cells.Item(51, 6).Value = 83
$ws3.Cells.Item(51, 7).Value = "L6173"

